$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Machine sheet: insert two new columns (x_location, y_location)
#    right before the existing column E ("from" / movement-related col),
#    pushing the old E..O block to G..Q.
# ---------------------------------------------------------------------
$machine = $wb.Worksheets.Item("Machine")

$machine.Range("E1:F1").EntireColumn.Insert()

# Match the width of the new columns to their left neighbour (column D),
# same as Excel does when inserting columns.
$machine.Columns("E:F").ColumnWidth = $machine.Range("D1").ColumnWidth

# Header row
$machine.Range("E1").Value = "x_location:int"
$machine.Range("F1").Value = "y_location:int"

# Data rows 2-13: x_location / y_location values (multiples of 5)
for ($row = 2; $row -le 13; $row++) {
    $val = ($row - 1) * 5
    $machine.Cells.Item($row, 5).Value = $val
    $machine.Cells.Item($row, 6).Value = $val
}

# ---------------------------------------------------------------------
# 2) Update the Machine sheet's _FilterDatabase defined name so its
#    range keeps up with the two newly inserted columns (O -> Q).
# ---------------------------------------------------------------------
$filterName = $wb.Names.Item("Machine!_FilterDatabase")
$filterName.RefersTo = "=Machine!`$A`$1:`$Q`$1"

# ---------------------------------------------------------------------
# 3) Selection bookkeeping: Location sheet gets a new selected cell
#    (E30) while Machine remains the active/tab-selected sheet with its
#    own new selection (E10). Doing Location first and Machine last
#    leaves Machine as the active tab, matching the saved workbook.
# ---------------------------------------------------------------------
$location = $wb.Worksheets.Item("Location")
$location.Range("E30").Select()

$machine.Range("E10").Select()
